$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 386, shifting existing rows 386:395 down to 387:396
$ws.Rows.Item(386).Insert()

# Populate the newly inserted row 386 with the new weekly price record
$ws.Cells.Item(386, 1).Value = 7
$ws.Cells.Item(386, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(386, 3).Value = "Ñuble"
$ws.Cells.Item(386, 4).Value = 45239
$ws.Cells.Item(386, 5).Value = 16
$ws.Cells.Item(386, 6).Value = 100112043
$ws.Cells.Item(386, 7).Value = "Pepino ensalada"
$ws.Cells.Item(386, 8).Value = "Sin especificar"
$ws.Cells.Item(386, 9).Value = "Primera"
$ws.Cells.Item(386, 10).Value = 70
$ws.Cells.Item(386, 11).Value = 14000
$ws.Cells.Item(386, 12).Value = 15000
$ws.Cells.Item(386, 13).Value = 14429
$ws.Cells.Item(386, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(386, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(386, 16).Value = 240
$ws.Cells.Item(386, 17).Value = 60
$ws.Cells.Item(386, 18).Value = "Hortaliza"
